# Fix the "Recorded By" (column G) values on the "Session Analysis Results"
# sheet: two specific comma-separated author lists need their tokens
# re-ordered. Every occurrence of the exact old text is replaced with the
# exact new text (find & replace by exact match), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of exact old "Recorded By" text -> exact new text.
$replacements = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "system, backup@backdoor.com, System" = "backup@backdoor.com, System, system"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row   # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Text

    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}

$wb.Save()
